$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the last data row (old row 11 / DF00002 X X entry) ---
# After this, the used range becomes A1:M10 (dimension shrinks automatically).
$ws.Rows(11).Delete()

# --- Row 2: refresh the computed "DATA_FIM_DT" value ---
$ws.Range("M2").Value = 45974
$ws.Range("M2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Row 6: populate the previously-empty "DATA_FIM_DT" value ---
$ws.Range("M6").Value = 45962
$ws.Range("M6").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Row 7: populate the previously-empty "DATA_FIM_DT" value ---
$ws.Range("M7").Value = 45962
$ws.Range("M7").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Row 8: new ticket data (was GO/TESTEGO, now DF/DIO1659) ---
$ws.Range("A8").Value = "DF"
$ws.Range("B8").Value = "DIO1659"
$ws.Range("H8").Value = "T - (Y 12/11/25_12H) - DF"
# I8 holds a date-like label that must stay plain text (not an auto-converted
# date serial) - force text entry with a quote prefix, then strip the
# resulting cell style back to Normal so no stray number format sticks.
$ws.Range("I8").Value = "'12/11/25"
$ws.Range("I8").Style = "Normal"
$ws.Range("M8").Value = 45973
$ws.Range("M8").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Row 9: new ticket data (was DF/DIO1659, now DF/DF00001) ---
$ws.Range("B9").Value = "DF00001"
$ws.Range("F9").Value = "TESTE1718"
$ws.Range("H9").Value = "T - (T 03/11/25_12H) - DF"
$ws.Range("I9").Value = "'03/11/25"
$ws.Range("I9").Style = "Normal"
$ws.Range("M9").Value = 45964
$ws.Range("M9").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Row 10: new ticket data (was DF/DF00001/TESTE1718, now DF/DF00002/X) ---
$ws.Range("B10").Value = "DF00002"
$ws.Range("F10").Value = "X"
$ws.Range("G10").Value = "X"
$ws.Range("H10").Value = "X - (X 03/11/25_12H) - DF"
$ws.Range("M10").Value = 45964
$ws.Range("M10").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output "SALDO_PECAS updated"
